$wb = $excel.ActiveWorkbook

# --- Step 1: Insert new "num" sheet before "perc" ---
$beforeSheet = $wb.Worksheets.Item("perc")
$numSheet = $wb.Worksheets.Add($beforeSheet)
$numSheet.Name = "num"

# --- Step 2: Update shared header/label strings across ALL sheets ---
$headerB = "FISM"
$headerC = "ITALIA"
$labels = @("tutte le pubblicazioni", "prevenzione primaria", "prevenzione secondaria", "prevenzione terziaria", "covid", "altre malattie")

# --- Sheet: num ---
$ws = $wb.Worksheets.Item("num")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 228
$ws.Cells.Item(2, 3).Value = 1400
$ws.Cells.Item(3, 2).Value = 25
$ws.Cells.Item(3, 3).Value = 132
$ws.Cells.Item(4, 2).Value = 41
$ws.Cells.Item(4, 3).Value = 172
$ws.Cells.Item(5, 2).Value = 47
$ws.Cells.Item(5, 3).Value = 191
$ws.Cells.Item(6, 2).Value = 25
$ws.Cells.Item(6, 3).Value = 117
$ws.Cells.Item(7, 2).Value = 5
$ws.Cells.Item(7, 3).Value = 37

# --- Sheet: perc ---
$ws = $wb.Worksheets.Item("perc")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 100
$ws.Cells.Item(2, 3).Value = 100
$ws.Cells.Item(3, 2).Value = 10.96491228070175
$ws.Cells.Item(3, 3).Value = 9.428571428571429
$ws.Cells.Item(4, 2).Value = 17.98245614035088
$ws.Cells.Item(4, 3).Value = 12.28571428571429
$ws.Cells.Item(5, 2).Value = 20.6140350877193
$ws.Cells.Item(5, 3).Value = 13.64285714285714
$ws.Cells.Item(6, 2).Value = 10.96491228070175
$ws.Cells.Item(6, 3).Value = 8.357142857142858
$ws.Cells.Item(7, 2).Value = 2.192982456140351
$ws.Cells.Item(7, 3).Value = 2.642857142857143

# --- Sheet: if ---
$ws = $wb.Worksheets.Item("if")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 8.419077981651375
$ws.Cells.Item(2, 3).Value = 6.149315499606608
$ws.Cells.Item(3, 2).Value = 10.64583333333333
$ws.Cells.Item(3, 3).Value = 6.506359374999999
$ws.Cells.Item(4, 2).Value = 13.95365853658537
$ws.Cells.Item(4, 3).Value = 9.149571428571429
$ws.Cells.Item(5, 2).Value = 6.106744680851063
$ws.Cells.Item(5, 3).Value = 4.963585635359116
$ws.Cells.Item(6, 2).Value = 6.056
$ws.Cells.Item(6, 3).Value = 5.221290598290597
$ws.Cells.Item(7, 2).Value = 5.475
$ws.Cells.Item(7, 3).Value = 6.597142857142857

# --- Sheet: altmetric ---
$ws = $wb.Worksheets.Item("altmetric")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 31.02631578947368
$ws.Cells.Item(2, 3).Value = 13.45746962115797
$ws.Cells.Item(3, 2).Value = 20.64
$ws.Cells.Item(3, 3).Value = 15.71969696969697
$ws.Cells.Item(4, 2).Value = 106.4390243902439
$ws.Cells.Item(4, 3).Value = 38.48255813953488
$ws.Cells.Item(5, 2).Value = 10.06382978723404
$ws.Cells.Item(5, 3).Value = 8.705263157894738
$ws.Cells.Item(6, 2).Value = 16.04
$ws.Cells.Item(6, 3).Value = 11.35042735042735
$ws.Cells.Item(7, 2).Value = 5.4
$ws.Cells.Item(7, 3).Value = 13.61111111111111

# --- Sheet: cima_index ---
$ws = $wb.Worksheets.Item("cima_index")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 2.352941176470588
$ws.Cells.Item(2, 3).Value = 2.086519114688129
$ws.Cells.Item(3, 2).Value = 2.047619047619047
$ws.Cells.Item(3, 3).Value = 1.99009900990099
$ws.Cells.Item(4, 2).Value = 2.966666666666667
$ws.Cells.Item(4, 3).Value = 2.384
$ws.Cells.Item(5, 2).Value = 2.696969696969697
$ws.Cells.Item(5, 3).Value = 2.13768115942029
$ws.Cells.Item(6, 2).Value = 3.05
$ws.Cells.Item(6, 3).Value = 2.293478260869565
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = 2

# --- Sheet: cit ---
$ws = $wb.Worksheets.Item("cit")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 14.47368421052632
$ws.Cells.Item(2, 3).Value = 10.56540385989993
$ws.Cells.Item(3, 2).Value = 21.08
$ws.Cells.Item(3, 3).Value = 14.79545454545454
$ws.Cells.Item(4, 2).Value = 12.78048780487805
$ws.Cells.Item(4, 3).Value = 17.47093023255814
$ws.Cells.Item(5, 2).Value = 11.78723404255319
$ws.Cells.Item(5, 3).Value = 9.726315789473684
$ws.Cells.Item(6, 2).Value = 41.68
$ws.Cells.Item(6, 3).Value = 21.90598290598291
$ws.Cells.Item(7, 2).Value = 11
$ws.Cells.Item(7, 3).Value = 13.86111111111111

# --- Sheet: reccit ---
$ws = $wb.Worksheets.Item("reccit")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 12.67105263157895
$ws.Cells.Item(2, 3).Value = 9.441029306647605
$ws.Cells.Item(3, 2).Value = 19.36
$ws.Cells.Item(3, 3).Value = 13.12878787878788
$ws.Cells.Item(4, 2).Value = 12
$ws.Cells.Item(4, 3).Value = 16.04651162790698
$ws.Cells.Item(5, 2).Value = 10.06382978723404
$ws.Cells.Item(5, 3).Value = 8.58421052631579
$ws.Cells.Item(6, 2).Value = 31.24
$ws.Cells.Item(6, 3).Value = 17.68376068376068
$ws.Cells.Item(7, 2).Value = 10
$ws.Cells.Item(7, 3).Value = 12.91666666666667

# --- Sheet: fcr ---
$ws = $wb.Worksheets.Item("fcr")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 7.931085714285715
$ws.Cells.Item(2, 3).Value = 6.806138996138996
$ws.Cells.Item(3, 2).Value = 10.2555
$ws.Cells.Item(3, 3).Value = 8.573076923076924
$ws.Cells.Item(4, 2).Value = 8.577272727272728
$ws.Cells.Item(4, 3).Value = 12.22378571428571
$ws.Cells.Item(5, 2).Value = 7.556111111111111
$ws.Cells.Item(5, 3).Value = 6.140645161290323
$ws.Cells.Item(6, 2).Value = 23.533
$ws.Cells.Item(6, 3).Value = 11.96726315789474
$ws.Cells.Item(7, 2).Value = 3.79
$ws.Cells.Item(7, 3).Value = 7.683225806451613

# --- Sheet: rcr ---
$ws = $wb.Worksheets.Item("rcr")
$ws.Range("B1").Value = $headerB
$ws.Range("C1").Value = $headerC
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}
$ws.Cells.Item(2, 2).Value = 2.387987804878049
$ws.Cells.Item(2, 3).Value = 2.076743455497382
$ws.Cells.Item(3, 2).Value = 2.6535
$ws.Cells.Item(3, 3).Value = 2.4452
$ws.Cells.Item(4, 2).Value = 2.246451612903226
$ws.Cells.Item(4, 3).Value = 3.197404580152672
$ws.Cells.Item(5, 2).Value = 2.9146875
$ws.Cells.Item(5, 3).Value = 2.283059701492538
$ws.Cells.Item(6, 2).Value = 5.193
$ws.Cells.Item(6, 3).Value = 3.262613636363636
$ws.Cells.Item(7, 2).Value = 1.63
$ws.Cells.Item(7, 3).Value = 2.553793103448276
